$d = $word.ActiveDocument

# Apply edits from the end of the document toward the start so that
# earlier (not-yet-processed) character offsets stay valid.

# -- 1. "Bitopi's" (curly apostrophe) -> "BBC's"  (offset 1360-1368)
$r = $d.Range(1360, 1368)
$r.Text = "BBC’s"

# -- 2. "Bitopi" -> "BBC"  (offset 1455-1461, "...at Bitopi.")
$r = $d.Range(1455, 1461)
$r.Text = "BBC"

# -- 3. "Software Engineer" -> "Software Developer"  (offset 1100-1117)
#    The run immediately after this one (a single space, 1117-1118) shares
#    identical run formatting, so a plain text replace would coalesce the
#    two runs together. Toggling Bold off/on/off around the edit keeps the
#    following run distinct, matching the original run layout.
$nbr = $d.Range(1117, 1118)
$nbr.Font.Bold = 1
$r = $d.Range(1100, 1117)
$r.Text = "Software Developer"
$nbr = $d.Range(1119, 1120)
$nbr.Font.Bold = 0

# -- 4. "Bitopi" -> "BBC"  (offset 541-547, "...apply my knowledge at Bitopi.")
$r = $d.Range(541, 547)
$r.Text = "BBC"

# -- 5. "Bitopi" -> "BBC"  (offset 308-314, "...position at Bitopi. With...")
$r = $d.Range(308, 314)
$r.Text = "BBC"

# -- 6. "Software Engineer" -> "Software Developer"  (offset 278-295)
#    Same neighbouring-run issue as edit 3 above.
$nbr = $d.Range(295, 296)
$nbr.Font.Bold = 1
$r = $d.Range(278, 295)
$r.Text = "Software Developer"
$nbr = $d.Range(297, 298)
$nbr.Font.Bold = 0

# -- 7. "Software Engineer" -> "Software Developer"  (offset 209-226, bold title)
#    This run is the last run in its paragraph, so no neighbour merge occurs.
$r = $d.Range(209, 226)
$r.Text = "Software Developer"

# -- 8. "26 February 2022" -> "27 February 2022"  (offset 176-192)
#    The preceding run ("Frankfurt am Main, ") has identical (empty) run
#    formatting, so toggle its Bold property off/on/off to keep it a
#    separate run instead of merging with the date text.
$nbr = $d.Range(157, 176)
$nbr.Font.Bold = 1
$r = $d.Range(176, 192)
$r.Text = "27 February 2022"
$nbr = $d.Range(157, 176)
$nbr.Font.Bold = 0

# -- 9. "Bitopi" -> "BBC"  (offset 132-138, standalone paragraph)
$r = $d.Range(132, 138)
$r.Text = "BBC"
